$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '58.770.15'),
    @('E2', '  +0.86%  '),
    @('D3', '2.493.08'),
    @('E3', '  +1.75%  '),
    @('D4', '1.00'),
    @('E4', '  -0.01%  '),
    @('D5', '534.46'),
    @('E5', '  +1.73%  '),
    @('D6', '136.32'),
    @('E6', '  +1.68%  '),
    @('E7', '  -0.17%  '),
    @('D8', '0.564'),
    @('E8', '  +2.24%  '),
    @('D9', '2.515.22'),
    @('E9', '  +2.21%  '),
    @('D10', '0.101'),
    @('E10', '  +2.80%  '),
    @('E11', '  -1.77%  '),
    @('D12', '5.38'),
    @('E12', '  +1.43%  '),
    @('E13', '  +2.28%  '),
    @('D14', '2.940.01'),
    @('E14', '  +1.74%  '),
    @('D15', '22.97'),
    @('E15', '  +2.34%  '),
    @('D16', '58.671.01'),
    @('E16', '  +0.86%  '),
    @('E17', '  +0.55%  '),
    @('D18', '2.519.27'),
    @('E18', '  +2.47%  '),
    @('D19', '11.10'),
    @('E19', '  +4.73%  '),
    @('D20', '4.25'),
    @('E20', '  +2.36%  '),
    @('D21', '322.50'),
    @('E21', '  +1.21%  '),
    @('D22', '1.00'),
    @('E22', '  +0.38%  '),
    @('E23', '  +5.14%  '),
    @('D24', '65.30'),
    @('E24', '  +5.23%  '),
    @('D25', '0.420'),
    @('E25', '  +4.03%  '),
    @('D26', '0.164'),
    @('E26', '  +0.85%  '),
    @('D27', '0.996'),
    @('E27', '  +1.46%  '),
    @('D28', '7.51'),
    @('E28', '  +1.11%  '),
    @('D29', '0.0₃0768'),
    @('E29', '  +2.95%  '),
    @('D30', '6.60'),
    @('E30', '  +2.41%  '),
    @('D31', '171.41'),
    @('E31', '  +5.17%  '),
    @('E32', '  +1.08%  '),
    @('E33', '  +10.54%  '),
    @('E34', '  -0.08%  '),
    @('B35', 'EthereumClassic'),
    @('C35', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'),
    @('D35', '18.34'),
    @('E35', '  +1.24%  '),
    @('B36', 'ImmutableX'),
    @('C36', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @('D36', '1.36'),
    @('E36', '  +1.44%  '),
    @('D37', '4.06'),
    @('E37', '  +1.95%  '),
    @('D38', '1.53'),
    @('E38', '  +0.31%  '),
    @('D39', '36.86'),
    @('E39', '  +1.38%  '),
    @('D40', '0.807'),
    @('E40', '  +3.91%  '),
    @('D41', '3.57'),
    @('E41', '  +1.88%  '),
    @('D42', '283.72'),
    @('E42', '  +5.18%  '),
    @('D43', '5.24'),
    @('E43', '  +5.24%  '),
    @('D44', '0.994'),
    @('E44', '  -0.54%  '),
    @('E45', '  +3.63%  '),
    @('D46', '129.47'),
    @('E46', '  +8.15%  '),
    @('E47', '  +0.27%  '),
    @('D48', '0.0921'),
    @('E48', '  +0.53%  '),
    @('E49', '  +0.59%  '),
    @('E50', '  +1.27%  '),
    @('D51', '17.33'),
    @('E51', '  +3.49%  ')
)

foreach ($u in $updates) {
    $cell = $ws.Range($u[0])
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.ClearFormats()
}
